$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1837349397590362
$ws.Range("C2").Value = 0.5963855421686747
$ws.Range("J2").Value = 0.006024096385542169
$ws.Range("P2").Value = 0.1174698795180723
$ws.Range("S2").Value = 0.0963855421686747
$ws.Range("B3").Value = 0.02898550724637681
$ws.Range("C3").Value = 0.04347826086956522
$ws.Range("J3").Value = 0.01932367149758454
$ws.Range("S3").Value = 0.2415458937198068
$ws.Range("J4").Value = 0.07407407407407407
$ws.Range("P4").Value = 0.5740740740740741
$ws.Range("S4").Value = 0.3518518518518519
$ws.Range("B6").Value = 0.08108108108108109
$ws.Range("D6").Value = 0.01351351351351351
$ws.Range("F6").Value = 0.06306306306306306
$ws.Range("J6").Value = 0.2522522522522522
$ws.Range("O6").Value = 0.02252252252252252
$ws.Range("Q6").Value = 0.1306306306306306
$ws.Range("R6").Value = 0.1261261261261261
$ws.Range("S6").Value = 0.3108108108108108
$ws.Range("B7").Value = 0.1100917431192661
$ws.Range("D7").Value = 0.02752293577981652
$ws.Range("F7").Value = 0.03211009174311927
$ws.Range("J7").Value = 0.1834862385321101
$ws.Range("O7").Value = 0.01834862385321101
$ws.Range("Q7").Value = 0.1926605504587156
$ws.Range("R7").Value = 0.07339449541284404
$ws.Range("S7").Value = 0.3623853211009174
$ws.Range("B8").Value = 0.108974358974359
$ws.Range("D8").Value = 0.02136752136752137
$ws.Range("E8").Value = 0.002136752136752137
$ws.Range("F8").Value = 0.05555555555555555
$ws.Range("J8").Value = 0.1581196581196581
$ws.Range("O8").Value = 0.01495726495726496
$ws.Range("Q8").Value = 0.1602564102564103
$ws.Range("R8").Value = 0.07264957264957266
$ws.Range("S8").Value = 0.405982905982906
$ws.Range("B9").Value = 0.06363636363636363
$ws.Range("D9").Value = 0.00909090909090909
$ws.Range("F9").Value = 0.04090909090909091
$ws.Range("J9").Value = 0.1636363636363636
$ws.Range("O9").Value = 0.02272727272727273
$ws.Range("Q9").Value = 0.15
$ws.Range("R9").Value = 0.07727272727272727
$ws.Range("S9").Value = 0.4727272727272727
$ws.Range("B10").Value = 0.1040424121935056
$ws.Range("D10").Value = 0.02120609675281643
$ws.Range("E10").Value = 0.0006626905235255136
$ws.Range("F10").Value = 0.06295559973492379
$ws.Range("J10").Value = 0.1491053677932405
$ws.Range("O10").Value = 0.01921802518223989
$ws.Range("Q10").Value = 0.1776010603048376
$ws.Range("R10").Value = 0.08151093439363817
$ws.Range("S10").Value = 0.3836978131212724
$ws.Range("G11").Value = 0.1359773371104816
$ws.Range("J11").Value = 0.1104815864022663
$ws.Range("K11").Value = 0.2067988668555241
$ws.Range("L11").Value = 0.5212464589235127
$ws.Range("S11").Value = 0.0254957507082153
$ws.Range("G12").Value = 0.7386934673366834
$ws.Range("J12").Value = 0.1758793969849246
$ws.Range("L12").Value = 0.02512562814070352
$ws.Range("S12").Value = 0.06030150753768844
$ws.Range("G13").Value = 0.6326530612244898
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.08163265306122448
$ws.Range("F15").Value = 0.01606425702811245
$ws.Range("H15").Value = 0.0963855421686747
$ws.Range("I15").Value = 0.09236947791164658
$ws.Range("J15").Value = 0.3734939759036144
$ws.Range("K15").Value = 0.04417670682730924
$ws.Range("M15").Value = 0.01204819277108434
$ws.Range("O15").Value = 0.07228915662650602
$ws.Range("S15").Value = 0.2931726907630522
$ws.Range("F16").Value = 0.02525252525252525
$ws.Range("H16").Value = 0.1616161616161616
$ws.Range("I16").Value = 0.07575757575757576
$ws.Range("J16").Value = 0.4141414141414141
$ws.Range("K16").Value = 0.1060606060606061
$ws.Range("M16").Value = 0.0101010101010101
$ws.Range("N16").Value = 0.005050505050505051
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.1464646464646465
$ws.Range("F17").Value = 0.006802721088435374
$ws.Range("H17").Value = 0.1904761904761905
$ws.Range("I17").Value = 0.09977324263038549
$ws.Range("J17").Value = 0.3968253968253968
$ws.Range("K17").Value = 0.1065759637188209
$ws.Range("M17").Value = 0.01133786848072562
$ws.Range("O17").Value = 0.04535147392290249
$ws.Range("S17").Value = 0.1428571428571428
$ws.Range("F18").Value = 0.02293577981651376
$ws.Range("H18").Value = 0.1605504587155963
$ws.Range("I18").Value = 0.1055045871559633
$ws.Range("J18").Value = 0.444954128440367
$ws.Range("K18").Value = 0.07798165137614679
$ws.Range("M18").Value = 0.01834862385321101
$ws.Range("N18").Value = 0.004587155963302753
$ws.Range("O18").Value = 0.06422018348623854
$ws.Range("S18").Value = 0.1009174311926606
$ws.Range("F19").Value = 0.01496598639455782
$ws.Range("H19").Value = 0.1986394557823129
$ws.Range("I19").Value = 0.07142857142857142
$ws.Range("J19").Value = 0.3795918367346939
$ws.Range("K19").Value = 0.117687074829932
$ws.Range("M19").Value = 0.02312925170068027
$ws.Range("O19").Value = 0.07823129251700681
$ws.Range("S19").Value = 0.1163265306122449
